$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "E2"  = 0.07975464681371225
    "C3"  = -3.017636378339217
    "E3"  = -3.229247082222797
    "C4"  = 0.003352386816724007
    "E4"  = -1.461031976610316
    "C5"  = 3.925837669383347
    "E5"  = 2.641604203902781
    "C6"  = 0.2381541440396262
    "E6"  = 1.60268309892857
    "C7"  = 4.993892964711621
    "E7"  = 2.260118192030736
    "C8"  = 6.711795724673664
    "E8"  = 6.409878804372982
    "C9"  = 0.5121603413743347
    "E9"  = 3.290935868252554
    "C10" = 1.745565778643887
    "E10" = 0.7985845180024986
    "C11" = 2.687500891103922
    "E11" = 1.922191950024699
    "C12" = 3.654655474034474
    "E12" = 3.068403604789749
    "C13" = 3.712036718632117
    "E13" = 3.908921577463587
    "C14" = 2.849400388885992
    "E14" = 3.535456592693387
    "C15" = -4.741003096464214
    "E15" = -2.156795995006056
    "C16" = 1.194925448553708
    "E16" = -2.068675356622807
    "C17" = -1.084365158506884
    "E17" = -1.089896342664354
    "C18" = -3.40787540386569
    "E18" = -2.191935020614488
    "C19" = -1.853660925652212
    "E19" = -1.878672029998096
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
